$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    5  = -4
    7  = -3
    11 = -5
    22 = -9
    23 = -1
    25 = -3
    28 = 1
    30 = 5
    32 = 1
    36 = -2
    39 = -5
    42 = -4
    46 = -6
    55 = 2
    60 = -2
    62 = 1
    63 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
